$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = "197640_ds"
$ws.Range("C9").Value  = "197640_us"
$ws.Range("C10").Value = "197658_us"
$ws.Range("C11").Value = "197658_ds2"
$ws.Range("C12").Value = "197658_ds1"
$ws.Range("C13").Value = "197662_us"
$ws.Range("C14").Value = "197662_ds"
$ws.Range("C15").Value = "197664_us"
$ws.Range("C16").Value = "197664_ds"
$ws.Range("C17").Value = "197665_ds"
$ws.Range("C25").Value = "197668_ds"
$ws.Range("C33").Value = "197663_us"
$ws.Range("C34").Value = "197663_ds"

$ws.Range("H46").Value = "``"
